# Auto-generated script applying scheduled market-data refresh values
# to the Adamantoise Profits workbook (one worksheet per crafting job).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 508
$ws.Range("I9").Value = 510
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 510
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = -341
$ws.Range("N9").Value = -838
$ws.Range("H15").Value = 2307.6562
$ws.Range("I15").Value = 2307.6562
$ws.Range("K15").Value = 6922.9686
$ws.Range("M15").Value = -6753.9686
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("H38").Value = 20228
$ws.Range("I38").Value = 285
$ws.Range("K38").Value = 855
$ws.Range("M38").Value = -483
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9752
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9142
$ws.Range("H123").Value = 84000
$ws.Range("J123").Value = 84000
$ws.Range("L123").Value = 84000
$ws.Range("N123").Value = -93800
$ws.Range("H138").Value = 2714.36
$ws.Range("I138").Value = 1016.9091
$ws.Range("J138").Value = 2924.1572
$ws.Range("K138").Value = 3050.7273
$ws.Range("L138").Value = 8772.471600000001
$ws.Range("M138").Value = 2089.2727
$ws.Range("N138").Value = -19052.4716

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4735.2
$ws.Range("I16").Value = 4419
$ws.Range("K16").Value = 4419
$ws.Range("M16").Value = -4132
$ws.Range("H44").Value = 68995.8
$ws.Range("J44").Value = 68994.75
$ws.Range("L44").Value = 68994.75
$ws.Range("N44").Value = -69970.75
$ws.Range("H122").Value = 3978.625
$ws.Range("I122").Value = 3460.093
$ws.Range("J122").Value = 5040.381
$ws.Range("K122").Value = 10380.279
$ws.Range("L122").Value = 15121.143
$ws.Range("M122").Value = -7930.278999999999
$ws.Range("N122").Value = -20021.143
$ws.Range("H123").Value = 62476.332
$ws.Range("J123").Value = 62476.332
$ws.Range("L123").Value = 62476.332
$ws.Range("N123").Value = -72276.33199999999
$ws.Range("H124").Value = 60369.2
$ws.Range("J124").Value = 60369.2
$ws.Range("L124").Value = 60369.2
$ws.Range("N124").Value = -70189.2
$ws.Range("H127").Value = 112205.43
$ws.Range("J127").Value = 112205.43
$ws.Range("L127").Value = 112205.43
$ws.Range("N127").Value = -122125.43
$ws.Range("H128").Value = 106804.664
$ws.Range("J128").Value = 106804.664
$ws.Range("L128").Value = 106804.664
$ws.Range("N128").Value = -116764.664
$ws.Range("H129").Value = 107455.8
$ws.Range("J129").Value = 107455.8
$ws.Range("L129").Value = 107455.8
$ws.Range("N129").Value = -117455.8
$ws.Range("H130").Value = 73806
$ws.Range("J130").Value = 73806
$ws.Range("L130").Value = 73806
$ws.Range("N130").Value = -83846

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("M55").Value = -9685
$ws.Range("H86").Value = 21492.945
$ws.Range("I86").Value = 19845.23
$ws.Range("J86").Value = 25777
$ws.Range("K86").Value = 19845.23
$ws.Range("L86").Value = 25777
$ws.Range("M86").Value = -18722.23
$ws.Range("N86").Value = -28023
$ws.Range("H87").Value = 71305.664
$ws.Range("J87").Value = 71305.664
$ws.Range("L87").Value = 71305.664
$ws.Range("N87").Value = -73677.664
$ws.Range("H89").Value = 21492.945
$ws.Range("I89").Value = 19845.23
$ws.Range("J89").Value = 25777
$ws.Range("K89").Value = 99226.14999999999
$ws.Range("L89").Value = 128885
$ws.Range("M89").Value = -93610.14999999999
$ws.Range("N89").Value = -140117
$ws.Range("H90").Value = 71305.664
$ws.Range("J90").Value = 71305.664
$ws.Range("L90").Value = 213916.992
$ws.Range("N90").Value = -225772.992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 432.33334
$ws.Range("I7").Value = 247.5
$ws.Range("K7").Value = 742.5
$ws.Range("M7").Value = -630.5
$ws.Range("H23").Value = 134.8
$ws.Range("I23").Value = 330
$ws.Range("J23").Value = 113.111115
$ws.Range("K23").Value = 990
$ws.Range("L23").Value = 339.333345
$ws.Range("M23").Value = -755
$ws.Range("N23").Value = -809.333345
$ws.Range("H107").Value = 658.8182
$ws.Range("J107").Value = 778.2
$ws.Range("L107").Value = 2334.6
$ws.Range("N107").Value = -6174.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 413.58334
$ws.Range("I2").Value = 377.5
$ws.Range("K2").Value = 377.5
$ws.Range("M2").Value = -264.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2631.6667
$ws.Range("I22").Value = 1822.5
$ws.Range("J22").Value = 4250
$ws.Range("K22").Value = 1822.5
$ws.Range("L22").Value = 4250
$ws.Range("M22").Value = -1527.5
$ws.Range("N22").Value = -4840
$ws.Range("H27").Value = 2631.6667
$ws.Range("I27").Value = 1822.5
$ws.Range("J27").Value = 4250
$ws.Range("K27").Value = 1822.5
$ws.Range("L27").Value = 4250
$ws.Range("M27").Value = -1715.5
$ws.Range("N27").Value = -4464
$ws.Range("H87").Value = 69994.5
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 69994.5
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 69994.5
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -72240.5
$ws.Range("H90").Value = 69994.5
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 69994.5
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 209983.5
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -221215.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 108466.336
$ws.Range("J16").Value = 108466.336
$ws.Range("L16").Value = 108466.336
$ws.Range("N16").Value = -109050.336
$ws.Range("H23").Value = 3352
$ws.Range("J23").Value = 4999
$ws.Range("L23").Value = 4999
$ws.Range("N23").Value = -5457
$ws.Range("H51").Value = 47733.5
$ws.Range("J51").Value = 59955
$ws.Range("L51").Value = 59955
$ws.Range("N51").Value = -60975
$ws.Range("H132").Value = 2817.647
$ws.Range("I132").Value = 2361.7
$ws.Range("J132").Value = 6237.25
$ws.Range("K132").Value = 7085.099999999999
$ws.Range("L132").Value = 18711.75
$ws.Range("M132").Value = -4555.099999999999
$ws.Range("N132").Value = -23771.75
